$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row: columns C and D swap (U and A)
$ws.Range("C1").Value = "U"
$ws.Range("D1").Value = "A"

# Update row 2 (FFR Lag): values change, columns follow new header order
$ws.Range("B2").Value = "0.29***"
$ws.Range("C2").Value = "-10.033***"
$ws.Range("D2").Value = "-0.171***"

# Row 3 becomes "U Lag" (was "A Lag")
$ws.Range("A3").Value = "U Lag"
$ws.Range("B3").Value = "0.025***"
$ws.Range("C3").Value = "-0.232***"
$ws.Range("D3").Value = "-0.027***"

# Row 4 becomes "A Lag" (was "U Lag")
$ws.Range("A4").Value = "A Lag"
$ws.Range("B4").Value = "0.365***"
$ws.Range("C4").Value = "-3.672***"
$ws.Range("D4").Value = "-0.667***"

# Remove rows 5 (Constant) and 6 (r2_adj) entirely
$ws.Range("A5:D6").EntireRow.Delete()
